# FINFLUX-2815  Stabilaizing automation script
#
# Updates the repayment numbers produced by a Make Repayment run (the
# overdue/variable-installment fee calc shifted by a cent or two) and adds
# the "Modify Transaction1" capture sheet that the automation adds after it
# edits (modifies) a transaction and navigates back to the loan.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - refreshed totals
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 1402.28
$summary.Range("E2").Value = 3597.72
$summary.Range("F2").Value = 276.18
$summary.Range("A3").Value = 181.66
$summary.Range("E3").Value = 84.68
$summary.Range("A5").Value = 0.74
$summary.Range("B5").Value = 0.74

# ---------------------------------------------------------------------
# 2. Repayment schedule sheet - refreshed totals
# ---------------------------------------------------------------------
$repay = $wb.Worksheets.Item("Repayment schedule")
$repay.Range("J4").Value = 0.74
$repay.Range("K4").Value = 888.46
$repay.Range("Q4").Value = 276.18
$repay.Range("F5").Value = 851.05
$repay.Range("G5").Value = 2470.4899999999998
$repay.Range("H5").Value = 36.67
$repay.Range("G6").Value = 1607.14
$repay.Range("G7").Value = 735.8
$repay.Range("F8").Value = 735.8
$repay.Range("K8").Value = 743.06
$repay.Range("Q8").Value = 743.06

# ---------------------------------------------------------------------
# 3. Transactions sheet - refreshed totals / ids
# ---------------------------------------------------------------------
$trans = $wb.Worksheets.Item("Transactions")
$trans.Range("A2").Value = 442
$trans.Range("F2").Value = 1402.28
$trans.Range("I2").Value = 0.74
$trans.Range("J2").Value = 3597.72
$trans.Range("A3").Value = 441

# ---------------------------------------------------------------------
# 4. New "Modify Transaction1" sheet, appended after "Transactions"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$modify = $wb.Worksheets.Add($null, $lastSheet)
$modify.Name = "Modify Transaction1"

$modify.Range("A1").Value = "OverDueTillDate"
$modify.Range("B1").Value = 42064
$modify.Range("A2").Value = "clickonsubmit"
$modify.Range("B2").Value = "Submit"
$modify.Range("A3").Value = "NavigateToLoan"
$modify.Range("B3").Value = "navigate"

# Label column - grey fill, matches the rest of the workbook's input labels
$labels = $modify.Range("A1:A3")
$labels.Interior.Color = 10921638
$labels.Font.Name = "Calibri"
$labels.Font.Size = 11

# Value column - green fill, Arial 10, matches the rest of the workbook
$values = $modify.Range("B1:B3")
$values.Interior.Color = 5296274
$values.Font.Name = "Calibri"
$values.Font.Size = 10

$modify.Range("B1").NumberFormat = "d-mmm-yy"
$modify.Range("B1").WrapText = $true

$modify.Columns.Item(1).ColumnWidth = 15.02
$modify.Columns.Item(2).ColumnWidth = 17.31

$modify.Range("D12").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Restore view state (selections / active tab) to match the source run
# ---------------------------------------------------------------------
$summary.Range("C8").Select() | Out-Null
$repay.Range("K7").Select() | Out-Null
$trans.Range("D8").Select() | Out-Null

$trans.Activate() | Out-Null
